$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are no longer part of the naive forecaster output
# (rows 2 and 3, column C; row 2, column E)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Update remaining values with corrected forecast figures (tiny
# floating-point corrections from the naive component forecaster bug fix)
$ws.Range("E3").Value = 4.422525088127305
$ws.Range("C4").Value = -14.45332333832744
$ws.Range("E4").Value = -2.928447329610051
$ws.Range("C5").Value = 8.600536527919612
$ws.Range("C7").Value = 4.639893381363192
$ws.Range("E8").Value = 2.429116709932599
$ws.Range("E9").Value = 3.941300050092877
$ws.Range("C12").Value = 4.695933104194361
$ws.Range("E15").Value = -1.352810423674367
$ws.Range("E16").Value = -0.3934198590721305
$ws.Range("C17").Value = 5.120680133083622
$ws.Range("C18").Value = -0.5532735011319123
$ws.Range("E19").Value = -1.28528149926006
